# Applies the "2022 column" update to the Лист1 sheet:
#  - adds a new data point for year 2022 in column S (mirroring column R's
#    layout/format for every row that already carries a 2021 figure)
#  - removes the placeholder blank S cells on the two section-header rows
#  - moves the active selection to R8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

function Copy-FormatAndSet {
    param(
        [string]$SourceAddr,
        [string]$TargetAddr,
        $Value
    )
    $ws.Range($SourceAddr).Copy() | Out-Null
    $ws.Range($TargetAddr).PasteSpecial($xlPasteFormats) | Out-Null
    if ($null -ne $Value) {
        $ws.Range($TargetAddr).Value2 = $Value
    }
}

# Header row: 2022 label, formatted like the other year headers (R4)
Copy-FormatAndSet "R4" "S4" 2022

# Country-wide total row (formatted like R5/R6)
Copy-FormatAndSet "R5" "S5" 33.152856050161155
Copy-FormatAndSet "R6" "S6" $null

# Data rows that already had an (empty) S placeholder -> fill them in
Copy-FormatAndSet "R7"  "S7"  32.831913512166025
Copy-FormatAndSet "R8"  "S8"  33.509346380994529
Copy-FormatAndSet "R10" "S10" 34.041194942162896
Copy-FormatAndSet "R11" "S11" 32.636018013483323
Copy-FormatAndSet "R18" "S18" 48.492370829119814
Copy-FormatAndSet "R19" "S19" 46.987664282528065
Copy-FormatAndSet "R20" "S20" 50.118899291215271
Copy-FormatAndSet "R21" "S21" 47.142900749295329
Copy-FormatAndSet "R22" "S22" 47.199946558584017

# Section-header rows 9 and 12 lose their blank S placeholder entirely
$ws.Range("S9").Clear() | Out-Null
$ws.Range("S12").Clear() | Out-Null

# Data rows that previously had no S cell at all -> add new ones
Copy-FormatAndSet "R13" "S13" 40.271414365477746
Copy-FormatAndSet "R14" "S14" 31.568157010024336
Copy-FormatAndSet "R15" "S15" 30.277813022272248
Copy-FormatAndSet "R16" "S16" 22.733608300917229

Copy-FormatAndSet "R23" "S23" 47.082025761639336
Copy-FormatAndSet "R24" "S24" 31.240016364696597
Copy-FormatAndSet "R25" "S25" 31.228685777194666
Copy-FormatAndSet "R26" "S26" 31.252112297543153
Copy-FormatAndSet "R27" "S27" 42.049857693482664
Copy-FormatAndSet "R28" "S28" 42.689244289315013
Copy-FormatAndSet "R29" "S29" 41.380596558931735
Copy-FormatAndSet "R30" "S30" 19.945481087558658
Copy-FormatAndSet "R31" "S31" 20.799187962023481
Copy-FormatAndSet "R32" "S32" 19.013188474520234
Copy-FormatAndSet "R33" "S33" 23.919779113642239
Copy-FormatAndSet "R34" "S34" 23.962040711070269
Copy-FormatAndSet "R35" "S35" 23.876854008981983
Copy-FormatAndSet "R36" "S36" 26.113584517813127
Copy-FormatAndSet "R37" "S37" 25.651528441631889
Copy-FormatAndSet "R38" "S38" 26.620973515499056
Copy-FormatAndSet "R39" "S39" 35.676666099583812
Copy-FormatAndSet "R40" "S40" 34.026766685280904
Copy-FormatAndSet "R41" "S41" 37.792274390474752
Copy-FormatAndSet "R42" "S42" 26.602385500795538
Copy-FormatAndSet "R43" "S43" 25.585637135242425

# Bottom total row (thick border bottom, like R44)
Copy-FormatAndSet "R44" "S44" 27.750206810614948

# Move the active selection the way it ended up in the authored workbook
$ws.Range("R8").Select() | Out-Null
